# Reorders the comma-separated "Recorded By" values in column G so that
# the list order is reversed (e.g. "dnasr281@gmail.com, System" becomes
# "System, dnasr281@gmail.com"), except for the specific combination
# "admin@admin.com, System" which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Text

    if ($value -eq $null -or $value -eq "") {
        continue
    }

    if ($value -eq "admin@admin.com, System") {
        continue
    }

    $parts = $value -split ", "

    if ($parts.Count -le 1) {
        continue
    }

    $n = $parts.Count
    $reversed = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $newValue = $reversed -join ", "

    $cell.Value = $newValue
}
